# removed ER tags from non-ER templates and non-ER tags
#
# The "SwateTemplateMetadata" sheet is renamed to "isa_template" and the
# ER (endpoint repository) entries (ER / ER Term Accession Number / ER Term
# Source REF values in B8:B10) are cleared, since this assay template is
# not tied to a specific endpoint repository.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Rename the metadata sheet.
$ws.Name = "isa_template"

# Clear the ER list values (ER / ER Term Accession Number / ER Term Source REF).
$ws.Range("B8:B10").ClearContents()

# Leave the metadata sheet active with B18 selected, matching the author's
# last recorded cursor position.
$ws.Activate()
$ws.Range("B18").Select() | Out-Null
